$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Straightforward numeric updates (country keeps its row position) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1103797
$ws.Range("C4").Value = 8774
$ws.Range("D4").Value = 156519
$ws.Range("E4").Value = 882826
$ws.Range("G4").Value = 596
$ws.Range("H4").Value = 64452

# Row 6: Italia
$ws.Range("B6").Value = 207428
$ws.Range("C6").Value = 1965
$ws.Range("D6").Value = 78249
$ws.Range("E6").Value = 100943
$ws.Range("F6").Value = 1578
$ws.Range("G6").Value = 269
$ws.Range("H6").Value = 28236

# Row 28: Singapur
$ws.Range("D28").Value = 1268
$ws.Range("E28").Value = 15817
$ws.Range("F28").Value = 23
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 16

# Row 45: Noruega
$ws.Range("B45").Value = 7770
$ws.Range("C45").Value = 32
$ws.Range("E45").Value = 7528

# Row 46: Chequia
$ws.Range("B46").Value = 7726
$ws.Range("C46").Value = 44
$ws.Range("D46").Value = 3359
$ws.Range("E46").Value = 4127
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 240

# Row 60: Luxemburgo
$ws.Range("B60").Value = 3802
$ws.Range("C60").Value = 18
$ws.Range("E60").Value = 497
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 92

# Row 103: Sri Lanka
$ws.Range("B103").Value = 690
$ws.Range("C103").Value = 27
$ws.Range("E103").Value = 521

# Row 116: Jordania
$ws.Range("B116").Value = 459
$ws.Range("C116").Value = 6
$ws.Range("D116").Value = 364
$ws.Range("E116").Value = 87

# --- Re-ranked pairs: new data pushes a country above its former neighbour ---

# Egipto overtakes Sudafrica (rows 52/53): Egipto now on row 52 with its
# updated figures, Sudafrica drops to row 53 keeping its previous figures.
$ws.Range("A52").Value = "Egipto"
$ws.Range("B52").Value = 5895
$ws.Range("C52").Value = 358
$ws.Range("D52").Value = 1460
$ws.Range("E52").Value = 4029
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 14
$ws.Range("H52").Value = 406

$ws.Range("A53").Value = "Sudafrica"
$ws.Range("B53").Value = 5647
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 2073
$ws.Range("E53").Value = 3471
$ws.Range("F53").Value = 36
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 103

# Zambia overtakes Suazilandia (rows 151/152): Zambia now on row 151 with
# its updated figures, Suazilandia drops to row 152 keeping its previous figures.
$ws.Range("A151").Value = "Zambia"
$ws.Range("B151").Value = 109
$ws.Range("C151").Value = 3
$ws.Range("D151").Value = 74
$ws.Range("E151").Value = 32
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 3

$ws.Range("A152").Value = "Suazilandia"
$ws.Range("B152").Value = 106
$ws.Range("C152").Value = 6
$ws.Range("D152").Value = 12
$ws.Range("E152").Value = 93
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 1

# --- Timestamp footer update ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 18:22"
